# Updates cryptos list values per the diff (price + volume(1h) columns,
# plus a few coin name/link/price swaps in rows 43/44/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.728.84"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "2.797.42"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.11%  "
$ws.Range("D15").Value = "3.240.05"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "2.802.75"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.958"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.96%  "
$ws.Range("D18").Value = "51.709.51"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("D22").Value = "0.0₃0972"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0909"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0454"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "121.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.56%  "
$ws.Range("E47").Value = "  +8.87%  "
$ws.Range("D48").Value = "2.120.57"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.978"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.84%  "
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.53%  "
